$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted price strings; force Text number format
# so Excel does not reinterpret numeric-looking strings as numbers/dates.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.927.56"
$ws.Range("E2").Value = "  +2.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.38"
$ws.Range("E3").Value = "  +2.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6389"
$ws.Range("E6").Value = "  +4.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2993"
$ws.Range("E8").Value = "  +4.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07484"
$ws.Range("E9").Value = "  +2.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.58"
$ws.Range("E10").Value = "  +7.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07680"
$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.885.07"
$ws.Range("E12").Value = "  +3.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.054"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6920"
$ws.Range("E14").Value = "  +5.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.36"
$ws.Range("E15").Value = "  +3.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009386"
$ws.Range("E16").Value = "  +4.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.069"
$ws.Range("E17").Value = "  +4.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.890.98"
$ws.Range("E18").Value = "  +2.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.128.52"
$ws.Range("E19").Value = "  +3.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.50"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.67"
$ws.Range("E21").Value = "  +2.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.380"
$ws.Range("E23").Value = "  +4.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.41"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1420"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.591"
$ws.Range("E27").Value = "  +2.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.98"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.505"
$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06068"
$ws.Range("E30").Value = "  +9.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.287"
$ws.Range("E31").Value = "  +7.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.136"
$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.147"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.890"
$ws.Range("E34").Value = "  +4.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.168"
$ws.Range("E35").Value = "  +3.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7283"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.608"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.854"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01796"
$ws.Range("E39").Value = "  +2.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.224.95"
$ws.Range("E40").Value = "  +1.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9307"
$ws.Range("E41").Value = "  +4.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.286"
$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.31"
$ws.Range("E46").Value = "  +2.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  +3.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5099"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.304"
$ws.Range("E49").Value = "  +3.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4096"
$ws.Range("E50").Value = "  +2.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1144"
$ws.Range("E51").Value = "  +3.36%  "

# Row 43 and 44 content swap (with updated price/volume)
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.036.82"
$ws.Range("E43").Value = "  +3.73%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.29%  "